$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (SM / Solar multiple) - revised nominal value and new bounds
$ws.Range("C5").Value = 2.4
$ws.Range("H5").Value = 1.4
$ws.Range("I5").Value = 3.8

# Row 15 (new STORAGE parameter: t_storage)
$ws.Range("B15").Value = "t_storage"
$ws.Range("C15").Value = 11
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 18

# Move the active selection on the bottom-right pane to I16 (matches author's cursor position)
$ws.Range("I16").Select()
